$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Student Name cell: "Mahesh Raj" -> "Janarthanan" + " K"
#
# In the target markup the name is split into two runs with Word's
# as-you-type proofing markup around the first word:
#
#   <w:proofErr w:type="spellStart"/>
#   <w:r><w:t>Janarthanan</w:t></w:r>
#   <w:proofErr w:type="spellEnd"/>
#   <w:r><w:t xml:space="preserve"> K</w:t></w:r>
#
# Find.Execute locates the paragraph holding "Mahesh Raj"; InsertXML then
# replaces that range with the exact OOXML above, inside a <w:p> that keeps
# the paragraph's original rsid attributes and TableParagraph style so only
# the run content (and the new proofErr bookmarks) actually change.
# ---------------------------------------------------------------------------
$nameRng = $d.Content
$foundName = $nameRng.Find.Execute("Mahesh Raj", $true, $false, $false, $false,
                                    $false, $true, 1, $false, "", 0)
if ($foundName) {
    $nameXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
               'w:rsidR="00E30DC0" w:rsidRDefault="00AD41DF">' +
               '<w:pPr><w:pStyle w:val="TableParagraph"/></w:pPr>' +
               '<w:proofErr w:type="spellStart"/>' +
               '<w:r><w:t>Janarthanan</w:t></w:r>' +
               '<w:proofErr w:type="spellEnd"/>' +
               '<w:r><w:t xml:space="preserve"> K</w:t></w:r>' +
               '</w:p>'
    $null = $nameRng.InsertXML($nameXml)
}

# ---------------------------------------------------------------------------
# 2. Student Roll Number cell: "952819106015" -> "952819106010"
# Plain text replace; run/character formatting (spacing, bookmark) is left
# untouched since only the digits inside <w:t> change.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("952819106015", $true, $false, $false, $false,
                                 $false, $true, 1, $false, "952819106010", 2)
